# Updated cryptos list with latest prices/volume percentages.
# Note: Price (column D) values are numeric-looking text (e.g. "70.055.01",
# "0.0000305") in the original sheet, so each is written with a leading
# apostrophe ('' in a single-quoted literal = one literal quote char) to force
# Excel to keep storing them as text instead of auto-converting to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''70.055.01'
$ws.Range('E2').Value = '  +0.24%  '
$ws.Range('D3').Value = '''3.552.50'
$ws.Range('E3').Value = '  +0.36%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '''603.48'
$ws.Range('E5').Value = '  -2.04%  '
$ws.Range('D6').Value = '''197.70'
$ws.Range('E6').Value = '  +6.84%  '
$ws.Range('E7').Value = '  -0.47%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '''0.210'
$ws.Range('E9').Value = '  -1.17%  '
$ws.Range('D10').Value = '''0.656'
$ws.Range('E10').Value = '  -0.12%  '
$ws.Range('D11').Value = '''54.27'
$ws.Range('E11').Value = '  +1.23%  '
$ws.Range('D12').Value = '''0.0000305'
$ws.Range('E12').Value = '  -0.83%  '
$ws.Range('D13').Value = '''9.58'
$ws.Range('E13').Value = '  +0.57%  '
$ws.Range('D14').Value = '''4.118.40'
$ws.Range('E14').Value = '  +0.24%  '
$ws.Range('D15').Value = '''602.32'
$ws.Range('E15').Value = '  -4.32%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '''70.207.39'
$ws.Range('E16').Value = '  +0.29%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').Value = '''19.17'
$ws.Range('E17').Value = '  +1.32%  '
$ws.Range('D18').Value = '''12.71'
$ws.Range('E18').Value = '  -1.11%  '
$ws.Range('D19').Value = '''3.557.84'
$ws.Range('E19').Value = '  -1.25%  '
$ws.Range('E21').Value = '  +0.42%  '
$ws.Range('D22').Value = '''18.02'
$ws.Range('E22').Value = '  +2.46%  '
$ws.Range('D23').Value = '''5.29'
$ws.Range('E23').Value = '  +8.01%  '
$ws.Range('D24').Value = '''103.37'
$ws.Range('E24').Value = '  +0.65%  '
$ws.Range('D25').Value = '''4.63'
$ws.Range('E25').Value = '  -2.41%  '
$ws.Range('D26').Value = '''3.14'
$ws.Range('E26').Value = '  +3.25%  '
$ws.Range('E27').Value = '  -0.19%  '
$ws.Range('D28').Value = '''9.69'
$ws.Range('E28').Value = '  +0.93%  '
$ws.Range('D29').Value = '''33.81'
$ws.Range('E29').Value = '  -1.60%  '
$ws.Range('D30').Value = '''4.54'
$ws.Range('E30').Value = '  +24.67%  '
$ws.Range('D31').Value = '''7.15'
$ws.Range('E31').Value = '  +1.58%  '
$ws.Range('D32').Value = '''12.72'
$ws.Range('E32').Value = '  +2.85%  '
$ws.Range('E33').Value = '  +0.87%  '
$ws.Range('D34').Value = '''63.44'
$ws.Range('E34').Value = '  -0.33%  '
$ws.Range('D35').Value = '''0.0₃0865'
$ws.Range('E35').Value = '  +11.34%  '
$ws.Range('D36').Value = '''3.757.02'
$ws.Range('E36').Value = '  +7.07%  '
$ws.Range('D37').Value = '''3.09'
$ws.Range('E37').Value = '  -4.01%  '
$ws.Range('E38').Value = '  +0.20%  '
$ws.Range('D39').Value = '''3.65'
$ws.Range('E39').Value = '  +3.43%  '
$ws.Range('D40').Value = '''0.396'
$ws.Range('E40').Value = '  -1.00%  '
$ws.Range('D41').Value = '''36.94'
$ws.Range('E41').Value = '  -0.37%  '
$ws.Range('D42').Value = '''490.23'
$ws.Range('E42').Value = '  -7.24%  '
$ws.Range('D43').Value = '''0.135'
$ws.Range('E43').Value = '  -1.12%  '
$ws.Range('E44').Value = '  -0.53%  '
$ws.Range('E45').Value = '  -3.23%  '
$ws.Range('E46').Value = '  -2.92%  '
$ws.Range('E47').Value = '  -1.35%  '
$ws.Range('D48').Value = '''1.00'
$ws.Range('E48').Value = '  +0.16%  '
$ws.Range('D49').Value = '''8.69'
$ws.Range('E49').Value = '  -4.09%  '
$ws.Range('D50').Value = '''0.000258'
$ws.Range('E50').Value = '  +7.25%  '
$ws.Range('D51').Value = '''1.32'
$ws.Range('E51').Value = '  +13.20%  '
